# Apply the weekly report update:
#  - Update the "Report Generated On" timestamp in D5
#  - Zero out the Total Billed Amount (C8) and per-line / total Pricing
#    figures (H16, H17) to reflect a no-violation / zero-billing scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report generation timestamp text in D5
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

# Zero out the billed amount / pricing figures
$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
